# Weekly data refresh: a new week's worth of price data (2 rows — "Primera"
# and "Segunda" quality grades) is inserted at row 169, pushing the existing
# rows 169:273 down to 171:275.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 169:170 (existing rows shift down).
$ws.Range("A169:R170").Insert()

# New row 169 — "Primera" quality.
$ws.Range("A169").Value = 7
$ws.Range("B169").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C169").Value = "Ñuble"
$ws.Range("D169").Value = 44806
$ws.Range("E169").Value = 16
$ws.Range("F169").Value = 100112009
$ws.Range("G169").Value = "Acelga"
$ws.Range("H169").Value = "Sin especificar"
$ws.Range("I169").Value = "Primera"
$ws.Range("J169").Value = 200
$ws.Range("K169").Value = 700
$ws.Range("L169").Value = 800
$ws.Range("M169").Value = 750
$ws.Range("N169").Value = "`$/atado 0,5 a 1 kilo"
$ws.Range("O169").Value = "Provincia de Diguillín"
$ws.Range("P169").Value = 750
$ws.Range("Q169").Value = 1
$ws.Range("R169").Value = "Hortaliza"

# New row 170 — "Segunda" quality.
$ws.Range("A170").Value = 7
$ws.Range("B170").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C170").Value = "Ñuble"
$ws.Range("D170").Value = 44806
$ws.Range("E170").Value = 16
$ws.Range("F170").Value = 100112009
$ws.Range("G170").Value = "Acelga"
$ws.Range("H170").Value = "Sin especificar"
$ws.Range("I170").Value = "Segunda"
$ws.Range("J170").Value = 150
$ws.Range("K170").Value = 600
$ws.Range("L170").Value = 600
$ws.Range("M170").Value = 600
$ws.Range("N170").Value = "`$/atado 0,5 a 1 kilo"
$ws.Range("O170").Value = "Provincia de Diguillín"
$ws.Range("P170").Value = 600
$ws.Range("Q170").Value = 1
$ws.Range("R170").Value = "Hortaliza"
